$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.087.00'
$ws.Range("E2").Value = '  +4.30%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.266.42'
$ws.Range("E3").Value = '  +2.82%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.61'
$ws.Range("E5").Value = '  +2.89%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.28'
$ws.Range("E6").Value = '  +3.26%  '

$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.605'
$ws.Range("E7").Value = '  -0.84%  '

$ws.Range("B8").Value = 'USDC'
$ws.Range("C8").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.263.78'

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.130'

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.75'
$ws.Range("E11").Value = '  +1.88%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.416'
$ws.Range("E12").Value = '  +4.75%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.838.54'
$ws.Range("E13").Value = '  +2.93%  '

$ws.Range("E14").Value = '  +0.85%  '

$ws.Range("E15").Value = '  +2.48%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '67.101.20'
$ws.Range("E16").Value = '  +4.34%  '

$ws.Range("E17").Value = '  +3.01%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.268.05'
$ws.Range("E18").Value = '  +2.80%  '

$ws.Range("E19").Value = '  +2.90%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.48'
$ws.Range("E20").Value = '  +2.86%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '372.10'
$ws.Range("E21").Value = '  +5.14%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.64'
$ws.Range("E22").Value = '  +5.91%  '

$ws.Range("E23").Value = '  -0.18%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.45'
$ws.Range("E24").Value = '  +3.33%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.514'
$ws.Range("E25").Value = '  +1.64%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.408.52'

$ws.Range("E28").Value = '  +1.31%  '

$ws.Range("E29").Value = '  +2.37%  '

$ws.Range("E30").Value = '  +0.11%  '

$ws.Range("E31").Value = '  +4.48%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.65'
$ws.Range("E32").Value = '  -0.69%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '22.66'
$ws.Range("E33").Value = '  +2.40%  '

$ws.Range("E34").Value = '  -0.05%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.27'
$ws.Range("E35").Value = '  +5.40%  '

$ws.Range("E36").Value = '  +2.38%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '168.13'
$ws.Range("E37").Value = '  +8.50%  '

$ws.Range("E38").Value = '  +4.43%  '

$ws.Range("E39").Value = '  +5.89%  '

$ws.Range("E40").Value = '  +9.75%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '27.34'
$ws.Range("E41").Value = '  +5.26%  '

$ws.Range("E42").Value = '  +3.05%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.763.55'
$ws.Range("E43").Value = '  +4.31%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.46'
$ws.Range("E44").Value = '  +6.88%  '

$ws.Range("E45").Value = '  +4.61%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '348.70'

$ws.Range("E47").Value = '  +3.36%  '

$ws.Range("E48").Value = '  +4.94%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '24.98'
$ws.Range("E49").Value = '  +5.03%  '

$ws.Range("E50").Value = '  +2.80%  '

$ws.Range("E51").Value = '  +2.49%  '
